$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers refreshed)
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687659435918"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168768942052"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168768943018"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168768990262"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687690606735"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687659076014.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687659256709.csv"
$ws1.Range("B4").Value = "go_stims-16511687659276812.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687659426234.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-16511687661777403.csv"
$ws2.Range("B3").Value = "OB-16511687672555246.csv"
$ws2.Range("B4").Value = "OB-16511687675015936.csv"
$ws2.Range("B5").Value = "ZB-match_5-16511687672119715.csv"
$ws2.Range("B6").Value = "OB-16511687677918491.csv"
$ws2.Range("B7").Value = "TB-1651168768919891.csv"
$ws2.Range("B8").Value = "ZB-match_2-165116876701092.csv"
$ws2.Range("B9").Value = "TB-16511687682204356.csv"
$ws2.Range("B10").Value = "TB-16511687681521716.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687689574463.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687689450176.csv"
$ws4.Range("B4").Value = "MM_stims-16511687689738166.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687689574463.csv"
$ws4.Range("B6").Value = "MM_stims-16511687689892979.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168768974821.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651168769030427.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687689941356.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687690050952.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511687690447617.csv"
